$wb = $excel.ActiveWorkbook

# --- Update the "Contact" sheet data ---
$wsContact = $wb.Worksheets.Item("Contact")
$wsContact.Range("A2").Value = "Activity Test External Contact"
$wsContact.Range("B2").Value = "ActivityCompany"

# --- Select range A2:B2 on the Contact sheet and make it the active sheet ---
$wsContact.Activate()
$wsContact.Range("A2:B2").Select()

# --- Restore selection on the MoreAttendees sheet (no longer the active tab) ---
$wsMoreAttendees = $wb.Worksheets.Item("MoreAttendees")
$wsMoreAttendees.Range("C8").Select()

# Re-activate Contact so it ends up as the active/selected tab in the saved file
$wsContact.Activate()
